# Update Mappings 22 Ontologies
# Adds a new "MS_DEF" column (F) after the existing MS_DESC column (E),
# with a literal "[]" value for every data row (2-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the bold/bordered style used by the other
# header cells in row 1 (copy style from E1).
$ws.Range("F1").Value = "MS_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

# Fill the new column for every existing data row with the literal "[]"
# string (matches column C / E cells: no special style).
for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 6).Value = "[]"
}
